$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "'64.032.15"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +0.31%  "
# Row 3
$ws.Range("D3").Value = "'2.761.27"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -0.04%  "
# Row 4
$ws.Range("E4").Value = "  +0.05%  "
# Row 5
$ws.Range("D5").Value = "'575.78"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -1.46%  "
# Row 6
$ws.Range("D6").Value = "'159.10"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.88%  "
# Row 7
$ws.Range("E7").Value = "  +0.10%  "
# Row 8
$ws.Range("D8").Value = "'0.601"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -2.84%  "
# Row 9
$ws.Range("E9").Value = "  -2.79%  "
# Row 10
$ws.Range("B10").Value = "Toncoin"
$ws.Range("C10").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D10").Value = "'5.84"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -13.61%  "
# Row 11
$ws.Range("B11").Value = "TRON"
$ws.Range("C11").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$ws.Range("D11").Value = "'0.165"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +3.53%  "
# Row 12
$ws.Range("E12").Value = "  -2.95%  "
# Row 13
$ws.Range("D13").Value = "'3.250.21"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -0.25%  "
# Row 14
$ws.Range("D14").Value = "'26.98"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -1.94%  "
# Row 15
$ws.Range("D15").Value = "'63.628.64"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -0.26%  "
# Row 16
$ws.Range("E16").Value = "  -4.88%  "
# Row 17
$ws.Range("D17").Value = "'2.767.44"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -0.25%  "
# Row 18
$ws.Range("D18").Value = "'12.11"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -2.10%  "
# Row 19
$ws.Range("D19").Value = "'4.84"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -3.01%  "
# Row 20
$ws.Range("D20").Value = "'359.84"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -1.62%  "
# Row 21
$ws.Range("E21").Value = "  -4.84%  "
# Row 22
$ws.Range("D22").Value = "'0.999"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.17%  "
# Row 23
$ws.Range("D23").Value = "'0.529"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -6.33%  "
# Row 24
$ws.Range("D24").Value = "'65.00"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -3.40%  "
# Row 25
$ws.Range("E25").Value = "  -3.66%  "
# Row 26
$ws.Range("B26").Value = "Binance-PegBSC-USD"
$ws.Range("C26").Value = "https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd"
$ws.Range("D26").Value = "'1.00"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +0.48%  "
# Row 27
$ws.Range("B27").Value = "InternetComputer(DFINITY)"
$ws.Range("C27").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D27").Value = "'8.52"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -1.91%  "
# Row 28
$ws.Range("D28").Value = "'0.0₃0905"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -5.69%  "
# Row 29
$ws.Range("D29").Value = "'7.35"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +1.12%  "
# Row 30
$ws.Range("E30").Value = "  -3.37%  "
# Row 31
$ws.Range("E31").Value = "  +0.10%  "
# Row 32
$ws.Range("D32").Value = "'170.54"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -1.28%  "
# Row 33
$ws.Range("B33").Value = "EthereumClassic"
$ws.Range("C33").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D33").Value = "'20.19"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -2.48%  "
# Row 34
$ws.Range("B34").Value = "NEARProtocol"
$ws.Range("C34").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D34").Value = "'4.91"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -2.55%  "
# Row 35
$ws.Range("B35").Value = "USDe"
$ws.Range("C35").Value = "https://coinranking.com/coin/exbfr2U-0+usde-usde"
$ws.Range("D35").Value = "'0.998"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +0.04%  "
# Row 36
$ws.Range("B36").Value = "ImmutableX"
$ws.Range("C36").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D36").Value = "'1.47"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -0.03%  "
# Row 37
$ws.Range("D37").Value = "'1.81"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -1.13%  "
# Row 38
$ws.Range("E38").Value = "  -1.81%  "
# Row 39
$ws.Range("D39").Value = "'348.36"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +2.76%  "
# Row 40
$ws.Range("D40").Value = "'6.25"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +0.50%  "
# Row 41
$ws.Range("E41").Value = "  -2.17%  "
# Row 42
$ws.Range("D42").Value = "'39.08"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -1.89%  "
# Row 43
$ws.Range("D43").Value = "'21.45"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -3.69%  "
# Row 44
$ws.Range("E44").Value = "  -3.83%  "
# Row 45
$ws.Range("E45").Value = "  -3.32%  "
# Row 46
$ws.Range("D46").Value = "'136.85"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -1.05%  "
# Row 47
$ws.Range("D47").Value = "'0.630"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -2.75%  "
# Row 48
$ws.Range("D48").Value = "'0.0253"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -2.99%  "
# Row 49
$ws.Range("E49").Value = "  -1.40%  "
# Row 50
$ws.Range("D50").Value = "'0.999"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +0.05%  "
# Row 51
$ws.Range("E51").Value = "  +0.21%  "
